$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Timetable")

# ---- Row 2: the CS161 / Dr. Sunil P V booking moves from X2:AC2 into AG2:AL2 ----
# (AG2:AL2 previously held the HS161 / Dr. Rajesh N S text; it gets overwritten.)
# X2:AC2 reverts to a block of plain, unmerged empty cells.
$ws.Range("X2:AC2").UnMerge()
$ws.Range("X2:AC2").ClearContents()
$ws.Range("X2:AC2").ClearFormats()

$ws.Range("AG2:AL2").Value = "CS161 | Problem Solving through Programming | Dr. Sunil P V | C002"

# ---- Row 5: part of the LUNCH BREAK block (R5:U5, plus the already-blank V5:W5)
#      becomes a class booking, taking on the look of the other CS161 slots.
#      Once merged, only the anchor cell (R5) should carry the value/format. ----
$ws.Range("R5:W5").Merge()
$ws.Range("R5:W5").Clear()
$ws.Range("R5:W5").Value = "CS161 | Problem Solving through Programming | Dr. Sunil C K | C002"

$ws.Range("I5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
